$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price and Volume(1h) columns)
# D-column (Price) values are forced to text to match the source data,
# which stores them as plain strings (e.g. "1.000", "27.142.14") rather
# than numbers -- Excel would otherwise auto-convert numeric-looking text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.142.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.718.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -6.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4699"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3434"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07274"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.041"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9997"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.856"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.718.65"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.888"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.94%  "
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06360"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("E21").Value = "  -3.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.627"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.180.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.60%  "
$ws.Range("E24").Value = "  -3.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.144"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.913.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.103"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.019"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09176"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.587"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.323"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02199"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05821"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.1994"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.739"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.62%  "
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5899"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.124"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.443"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5647"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.15%  "
$ws.Range("E46").Value = "  -5.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.838"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06647"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.086"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("E51").Value = "  -0.07%  "
